{"js": "// Change the Minedu (Ministry of Education) name reference in the body text:\n// \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2, \u0388\u03c1\u03b5\u03c5\u03bd\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\" -> \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\"\n// i.e. drop the \", \u0388\u03c1\u03b5\u03c5\u03bd\u03b1\u03c2\" (\", Research\") part of the ministry's old name,\n// inside the \"\u039c\u03b5\u03c4\u03ac \u03c4\u03b7\u03bd \u03b5\u03c0\u03b9\u03c3\u03c4\u03c1\u03bf\u03c6\u03ae ...\" paragraph of the school transports template.\n\nconst oldText = \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2, \u0388\u03c1\u03b5\u03c5\u03bd\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\";\nconst newText = \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + oldText);\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Change the Minedu (Ministry of Education) name reference in the body text:\n# \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2, \u0388\u03c1\u03b5\u03c5\u03bd\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\" -> \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\"\n# i.e. drop the \", \u0388\u03c1\u03b5\u03c5\u03bd\u03b1\u03c2\" (\", Research\") part of the ministry's old name.\n\n$d = $word.ActiveDocument\n\n$oldText = \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2, \u0388\u03c1\u03b5\u03c5\u03bd\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\"\n$newText = \"\u03a0\u03b1\u03b9\u03b4\u03b5\u03af\u03b1\u03c2 \u03ba\u03b1\u03b9 \u0398\u03c1\u03b7\u03c3\u03ba\u03b5\u03c5\u03bc\u03ac\u03c4\u03c9\u03bd\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n"}
